$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("D1").Value = "Moyenne Factuality"
$ws.Range("E1").Value = "Moyenne Readability"

# Update row 2 contents
$ws.Range("B2").Value = "Factuality and Readability"
$ws.Range("C2").Value = "0.1585407853126526 / 2"
$ws.Range("D2").Value = 0.1585407853126526
$ws.Range("E2").Value = 2

# Remove the now-obsolete row 3 entirely
$ws.Rows.Item(3).Delete()
